# Natmi following Dr Hou advice
# Recompute the Wnt4-Fzd2 LR-pair sheet to include the full 4x4
# sending-cluster x target-cluster grid (ECs, FAPs, M2, sCs),
# updating existing rows 2-13 in place and appending new rows 14-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Wnt4"
$ws.Cells.Item(2,3).Value = "Fzd2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.437774666666667
$ws.Cells.Item(2,8).Value = 4.313324
$ws.Cells.Item(2,9).Value = 0.2965885866872326
$ws.Cells.Item(2,10).Value = 0.2965885866872326
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.2081856666666667
$ws.Cells.Item(2,14).Value = 0.624557
$ws.Cells.Item(2,15).Value = 0.0127208067884984
$ws.Cells.Item(2,16).Value = 0.0127208067884984
$ws.Cells.Item(2,17).Value = 0.2993240774964445
$ws.Cells.Item(2,18).Value = 2.693916697468
$ws.Cells.Item(2,19).Value = 0.003772846106922096
$ws.Cells.Item(2,20).Value = 0.003772846106922096

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Wnt4"
$ws.Cells.Item(3,3).Value = "Fzd2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.437774666666667
$ws.Cells.Item(3,8).Value = 4.313324
$ws.Cells.Item(3,9).Value = 0.2965885866872326
$ws.Cells.Item(3,10).Value = 0.2965885866872326
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 13.42533133333333
$ws.Cells.Item(3,14).Value = 40.275994
$ws.Cells.Item(3,15).Value = 0.8203304708596988
$ws.Cells.Item(3,16).Value = 0.8203304708596987
$ws.Cells.Item(3,17).Value = 19.30260128267289
$ws.Cells.Item(3,18).Value = 173.723411544056
$ws.Cells.Item(3,19).Value = 0.2433006549687501
$ws.Cells.Item(3,20).Value = 0.2433006549687501

# Row 4: ECs -> M2
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Wnt4"
$ws.Cells.Item(4,3).Value = "Fzd2"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.437774666666667
$ws.Cells.Item(4,8).Value = 4.313324
$ws.Cells.Item(4,9).Value = 0.2965885866872326
$ws.Cells.Item(4,10).Value = 0.2965885866872326
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.04883999999999999
$ws.Cells.Item(4,14).Value = 0.14652
$ws.Cells.Item(4,15).Value = 0.002984279434304292
$ws.Cells.Item(4,16).Value = 0.002984279434304292
$ws.Cells.Item(4,17).Value = 0.07022091471999999
$ws.Cells.Item(4,18).Value = 0.6319882324799999
$ws.Cells.Item(4,19).Value = 0.0008851032197000839
$ws.Cells.Item(4,20).Value = 0.000885103219700084

# Row 5: ECs -> sCs
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Wnt4"
$ws.Cells.Item(5,3).Value = "Fzd2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.437774666666667
$ws.Cells.Item(5,8).Value = 4.313324
$ws.Cells.Item(5,9).Value = 0.2965885866872326
$ws.Cells.Item(5,10).Value = 0.2965885866872326
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 2.683402666666666
$ws.Cells.Item(5,14).Value = 8.050208
$ws.Cells.Item(5,15).Value = 0.1639644429174985
$ws.Cells.Item(5,16).Value = 0.1639644429174985
$ws.Cells.Item(5,17).Value = 3.858128374599111
$ws.Cells.Item(5,18).Value = 34.723155371392
$ws.Cells.Item(5,19).Value = 0.04862998239186032
$ws.Cells.Item(5,20).Value = 0.04862998239186032

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Wnt4"
$ws.Cells.Item(6,3).Value = "Fzd2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.589574
$ws.Cells.Item(6,8).Value = 4.768721999999999
$ws.Cells.Item(6,9).Value = 0.3279022207198701
$ws.Cells.Item(6,10).Value = 0.3279022207198702
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.2081856666666667
$ws.Cells.Item(6,14).Value = 0.624557
$ws.Cells.Item(6,15).Value = 0.0127208067884984
$ws.Cells.Item(6,16).Value = 0.0127208067884984
$ws.Cells.Item(6,17).Value = 0.330926522906
$ws.Cells.Item(6,18).Value = 2.978338706154
$ws.Cells.Item(6,19).Value = 0.004171180795297025
$ws.Cells.Item(6,20).Value = 0.004171180795297026

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Wnt4"
$ws.Cells.Item(7,3).Value = "Fzd2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.589574
$ws.Cells.Item(7,8).Value = 4.768721999999999
$ws.Cells.Item(7,9).Value = 0.3279022207198701
$ws.Cells.Item(7,10).Value = 0.3279022207198702
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 13.42533133333333
$ws.Cells.Item(7,14).Value = 40.275994
$ws.Cells.Item(7,15).Value = 0.8203304708596988
$ws.Cells.Item(7,16).Value = 0.8203304708596987
$ws.Cells.Item(7,17).Value = 21.340557628852
$ws.Cells.Item(7,18).Value = 192.065018659668
$ws.Cells.Item(7,19).Value = 0.268988183119072
$ws.Cells.Item(7,20).Value = 0.268988183119072

# Row 8: FAPs -> M2
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Wnt4"
$ws.Cells.Item(8,3).Value = "Fzd2"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.589574
$ws.Cells.Item(8,8).Value = 4.768721999999999
$ws.Cells.Item(8,9).Value = 0.3279022207198701
$ws.Cells.Item(8,10).Value = 0.3279022207198702
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.04883999999999999
$ws.Cells.Item(8,14).Value = 0.14652
$ws.Cells.Item(8,15).Value = 0.002984279434304292
$ws.Cells.Item(8,16).Value = 0.002984279434304292
$ws.Cells.Item(8,17).Value = 0.07763479415999998
$ws.Cells.Item(8,18).Value = 0.6987131474399998
$ws.Cells.Item(8,19).Value = 0.0009785518537570152
$ws.Cells.Item(8,20).Value = 0.0009785518537570152

# Row 9: FAPs -> sCs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Wnt4"
$ws.Cells.Item(9,3).Value = "Fzd2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.589574
$ws.Cells.Item(9,8).Value = 4.768721999999999
$ws.Cells.Item(9,9).Value = 0.3279022207198701
$ws.Cells.Item(9,10).Value = 0.3279022207198702
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.683402666666666
$ws.Cells.Item(9,14).Value = 8.050208
$ws.Cells.Item(9,15).Value = 0.1639644429174985
$ws.Cells.Item(9,16).Value = 0.1639644429174985
$ws.Cells.Item(9,17).Value = 4.265467110463999
$ws.Cells.Item(9,18).Value = 38.38920399417599
$ws.Cells.Item(9,19).Value = 0.05376430495174415
$ws.Cells.Item(9,20).Value = 0.05376430495174416

# Row 10: M2 -> ECs
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Wnt4"
$ws.Cells.Item(10,3).Value = "Fzd2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.104440333333333
$ws.Cells.Item(10,8).Value = 3.313321
$ws.Cells.Item(10,9).Value = 0.2278273537140099
$ws.Cells.Item(10,10).Value = 0.22782735371401
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.2081856666666667
$ws.Cells.Item(10,14).Value = 0.624557
$ws.Cells.Item(10,15).Value = 0.0127208067884984
$ws.Cells.Item(10,16).Value = 0.0127208067884984
$ws.Cells.Item(10,17).Value = 0.2299286470885556
$ws.Cells.Item(10,18).Value = 2.069357823797
$ws.Cells.Item(10,19).Value = 0.002898147747730805
$ws.Cells.Item(10,20).Value = 0.002898147747730805

# Row 11: M2 -> FAPs
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Wnt4"
$ws.Cells.Item(11,3).Value = "Fzd2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.104440333333333
$ws.Cells.Item(11,8).Value = 3.313321
$ws.Cells.Item(11,9).Value = 0.2278273537140099
$ws.Cells.Item(11,10).Value = 0.22782735371401
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 13.42533133333333
$ws.Cells.Item(11,14).Value = 40.275994
$ws.Cells.Item(11,15).Value = 0.8203304708596988
$ws.Cells.Item(11,16).Value = 0.8203304708596987
$ws.Cells.Item(11,17).Value = 14.82747741289711
$ws.Cells.Item(11,18).Value = 133.447296716074
$ws.Cells.Item(11,19).Value = 0.1868937203469329
$ws.Cells.Item(11,20).Value = 0.1868937203469329

# Row 12: M2 -> M2
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Wnt4"
$ws.Cells.Item(12,3).Value = "Fzd2"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.104440333333333
$ws.Cells.Item(12,8).Value = 3.313321
$ws.Cells.Item(12,9).Value = 0.2278273537140099
$ws.Cells.Item(12,10).Value = 0.22782735371401
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.04883999999999999
$ws.Cells.Item(12,14).Value = 0.14652
$ws.Cells.Item(12,15).Value = 0.002984279434304292
$ws.Cells.Item(12,16).Value = 0.002984279434304292
$ws.Cells.Item(12,17).Value = 0.05394086588
$ws.Cells.Item(12,18).Value = 0.48546779292
$ws.Cells.Item(12,19).Value = 0.0006799004862606895
$ws.Cells.Item(12,20).Value = 0.0006799004862606895

# Row 13: M2 -> sCs
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Wnt4"
$ws.Cells.Item(13,3).Value = "Fzd2"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.104440333333333
$ws.Cells.Item(13,8).Value = 3.313321
$ws.Cells.Item(13,9).Value = 0.2278273537140099
$ws.Cells.Item(13,10).Value = 0.22782735371401
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 2.683402666666666
$ws.Cells.Item(13,14).Value = 8.050208
$ws.Cells.Item(13,15).Value = 0.1639644429174985
$ws.Cells.Item(13,16).Value = 0.1639644429174985
$ws.Cells.Item(13,17).Value = 2.963658135640889
$ws.Cells.Item(13,18).Value = 26.672923220768
$ws.Cells.Item(13,19).Value = 0.03735558513308553
$ws.Cells.Item(13,20).Value = 0.03735558513308554

# Row 14: sCs -> ECs
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Wnt4"
$ws.Cells.Item(14,3).Value = "Fzd2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.7159183333333333
$ws.Cells.Item(14,8).Value = 2.147755
$ws.Cells.Item(14,9).Value = 0.1476818388788872
$ws.Cells.Item(14,10).Value = 0.1476818388788872
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.2081856666666667
$ws.Cells.Item(14,14).Value = 0.624557
$ws.Cells.Item(14,15).Value = 0.0127208067884984
$ws.Cells.Item(14,16).Value = 0.0127208067884984
$ws.Cells.Item(14,17).Value = 0.1490439355038889
$ws.Cells.Item(14,18).Value = 1.341395419535
$ws.Cells.Item(14,19).Value = 0.001878632138548476
$ws.Cells.Item(14,20).Value = 0.001878632138548476

# Row 15: sCs -> FAPs
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Wnt4"
$ws.Cells.Item(15,3).Value = "Fzd2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.7159183333333333
$ws.Cells.Item(15,8).Value = 2.147755
$ws.Cells.Item(15,9).Value = 0.1476818388788872
$ws.Cells.Item(15,10).Value = 0.1476818388788872
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 13.42533133333333
$ws.Cells.Item(15,14).Value = 40.275994
$ws.Cells.Item(15,15).Value = 0.8203304708596988
$ws.Cells.Item(15,16).Value = 0.8203304708596987
$ws.Cells.Item(15,17).Value = 9.611440832607776
$ws.Cells.Item(15,18).Value = 86.50296749347
$ws.Cells.Item(15,19).Value = 0.1211479124249437
$ws.Cells.Item(15,20).Value = 0.1211479124249437

# Row 16: sCs -> M2
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Wnt4"
$ws.Cells.Item(16,3).Value = "Fzd2"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.7159183333333333
$ws.Cells.Item(16,8).Value = 2.147755
$ws.Cells.Item(16,9).Value = 0.1476818388788872
$ws.Cells.Item(16,10).Value = 0.1476818388788872
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.04883999999999999
$ws.Cells.Item(16,14).Value = 0.14652
$ws.Cells.Item(16,15).Value = 0.002984279434304292
$ws.Cells.Item(16,16).Value = 0.002984279434304292
$ws.Cells.Item(16,17).Value = 0.03496545139999999
$ws.Cells.Item(16,18).Value = 0.3146890626
$ws.Cells.Item(16,19).Value = 0.0004407238745865031
$ws.Cells.Item(16,20).Value = 0.0004407238745865031

# Row 17: sCs -> sCs
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Wnt4"
$ws.Cells.Item(17,3).Value = "Fzd2"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.7159183333333333
$ws.Cells.Item(17,8).Value = 2.147755
$ws.Cells.Item(17,9).Value = 0.1476818388788872
$ws.Cells.Item(17,10).Value = 0.1476818388788872
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 2.683402666666666
$ws.Cells.Item(17,14).Value = 8.050208
$ws.Cells.Item(17,15).Value = 0.1639644429174985
$ws.Cells.Item(17,16).Value = 0.1639644429174985
$ws.Cells.Item(17,17).Value = 1.921097164782222
$ws.Cells.Item(17,18).Value = 17.28987448304
$ws.Cells.Item(17,19).Value = 0.02421457044080852
$ws.Cells.Item(17,20).Value = 0.02421457044080852

